$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Notes master "date" field placeholder: 02/01/2021 -> 2/1/21
# ------------------------------------------------------------------
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$hf.DateAndTime.Text = "2/1/21"

# ------------------------------------------------------------------
# 2) Notes master "slide number" field placeholder: <Nr.> -> <#>
# ------------------------------------------------------------------
$hf.SlideNumber.Text = [char]0x2039 + "#" + [char]0x203A

# ------------------------------------------------------------------
# 3) Fig. 1 caption: split into two runs and add a trailing period
#    "Fig. 1 Search terms visualized by icons"
#    -> "Fig. 1 Search terms visualized " + "by icons." (italic run)
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(12)
$tr = $sh.TextFrame.TextRange

if ($tr.Text -eq "Fig. 1 Search terms visualized by icons") {
    $sub = $tr.Characters(32, 8)
    $sub.Text = "by icons."
}
